# Update 'F' column (想去人数 / interested-count) values across all four sheets
# of the workbook, per the commit 'Update gh-pages to output generated at 456a3b4'.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 7791  # was 7786
$ws.Cells.Item(3, 6).Value = 7791  # was 7786
$ws.Cells.Item(5, 6).Value = 7950  # was 7944
$ws.Cells.Item(6, 6).Value = 45  # was 44
$ws.Cells.Item(9, 6).Value = 6828  # was 6819
$ws.Cells.Item(10, 6).Value = 3420  # was 3419
$ws.Cells.Item(20, 6).Value = 68  # was 66
$ws.Cells.Item(23, 6).Value = 2  # was 1
$ws.Cells.Item(25, 6).Value = 3926  # was 3923
$ws.Cells.Item(28, 6).Value = 1015  # was 1014
$ws.Cells.Item(29, 6).Value = 503  # was 501
$ws.Cells.Item(30, 6).Value = 1541  # was 1538
$ws.Cells.Item(32, 6).Value = 71  # was 70
$ws.Cells.Item(33, 6).Value = 2829  # was 2827
$ws.Cells.Item(34, 6).Value = 1993  # was 1987
$ws.Cells.Item(36, 6).Value = 60  # was 59
$ws.Cells.Item(38, 6).Value = 101  # was 100
$ws.Cells.Item(39, 6).Value = 3848  # was 3840
$ws.Cells.Item(40, 6).Value = 349  # was 347
$ws.Cells.Item(45, 6).Value = 37  # was 25
$ws.Cells.Item(46, 6).Value = 1482  # was 1477

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(7, 6).Value = 47  # was 46
$ws.Cells.Item(17, 6).Value = 267  # was 248

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 144  # was 143

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 144  # was 143
$ws.Cells.Item(5, 6).Value = 7791  # was 7786
$ws.Cells.Item(6, 6).Value = 7791  # was 7786
$ws.Cells.Item(7, 6).Value = 7950  # was 7944
$ws.Cells.Item(8, 6).Value = 45  # was 44
$ws.Cells.Item(11, 6).Value = 6828  # was 6819
$ws.Cells.Item(12, 6).Value = 3420  # was 3419
$ws.Cells.Item(19, 6).Value = 47  # was 46
$ws.Cells.Item(20, 6).Value = 68  # was 66
$ws.Cells.Item(24, 6).Value = 2  # was 1
$ws.Cells.Item(26, 6).Value = 3926  # was 3923
$ws.Cells.Item(31, 6).Value = 503  # was 501
$ws.Cells.Item(32, 6).Value = 1541  # was 1538
$ws.Cells.Item(34, 6).Value = 71  # was 70
$ws.Cells.Item(35, 6).Value = 1993  # was 1987
$ws.Cells.Item(37, 6).Value = 60  # was 59
$ws.Cells.Item(38, 6).Value = 101  # was 100
$ws.Cells.Item(40, 6).Value = 3848  # was 3840
$ws.Cells.Item(41, 6).Value = 349  # was 347
$ws.Cells.Item(46, 6).Value = 37  # was 25
$ws.Cells.Item(47, 6).Value = 1482  # was 1477
